# Scheduled runner update: refresh leve-profit market price figures across
# all class sheets (currentAveragePrice* / LevePrice* / LeveProfit* cols).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2
$ws.Range("J9").Value = 2
$ws.Range("L9").Value = 2
$ws.Range("N9").Value = -340

$ws.Range("H19").Value = 4888.375
$ws.Range("I19").Value = 4875.5
$ws.Range("K19").Value = 4875.5
$ws.Range("M19").Value = -4700.5

$ws.Range("H99").Value = 704.3570999999999
$ws.Range("I99").Value = 450.84616
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 1352.53848
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = 145.4615200000001
$ws.Range("N99").Value = -14996

$ws.Range("H101").Value = 812.4286
$ws.Range("I101").Value = 625.4
$ws.Range("J101").Value = 1280
$ws.Range("K101").Value = 1876.2
$ws.Range("L101").Value = 3840
$ws.Range("M101").Value = -254.1999999999998
$ws.Range("N101").Value = -7084

$ws.Range("H112").Value = 1220.3793
$ws.Range("J112").Value = 1245.8462
$ws.Range("L112").Value = 3737.5386
$ws.Range("N112").Value = -5953.5386

$ws.Range("H132").Value = 50799.332
$ws.Range("I132").Value = 79207.08
$ws.Range("J132").Value = 4636.75
$ws.Range("K132").Value = 237621.24
$ws.Range("L132").Value = 13910.25
$ws.Range("M132").Value = -235091.24
$ws.Range("N132").Value = -18970.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1783.0294
$ws.Range("I32").Value = 1812.9697
$ws.Range("J32").Value = 795
$ws.Range("K32").Value = 1812.9697
$ws.Range("L32").Value = 795
$ws.Range("M32").Value = -1525.9697
$ws.Range("N32").Value = -1369

$ws.Range("H45").Value = 1996.1666
$ws.Range("I45").Value = 1996.1666
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1996.1666
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -1619.1666

$ws.Range("H61").Value = 2078.6667
$ws.Range("I61").Value = 1994.5
$ws.Range("K61").Value = 1994.5
$ws.Range("M61").Value = -1782.5

$ws.Range("H97").Value = 493.06668
$ws.Range("I97").Value = 540.1667
$ws.Range("K97").Value = 540.1667
$ws.Range("M97").Value = -44.16669999999999

$ws.Range("H132").Value = 15156558
$ws.Range("I132").Value = 3039.3914
$ws.Range("K132").Value = 9118.174199999999
$ws.Range("M132").Value = -6588.174199999999

$ws.Range("H136").Value = 2078.6667
$ws.Range("I136").Value = 1994.5
$ws.Range("K136").Value = 5983.5
$ws.Range("M136").Value = -3433.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1807.7646
$ws.Range("J20").Value = 2554
$ws.Range("L20").Value = 2554
$ws.Range("N20").Value = -3048

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H99").Value = 2262.7856
$ws.Range("I99").Value = 2302.3076
$ws.Range("J99").Value = 1749
$ws.Range("K99").Value = 2302.3076
$ws.Range("L99").Value = 1749
$ws.Range("M99").Value = -804.3076000000001
$ws.Range("N99").Value = -4745

$ws.Range("H134").Value = 26522478
$ws.Range("I134").Value = 13165977
$ws.Range("K134").Value = 39497931
$ws.Range("M134").Value = -39495396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1750
$ws.Range("J16").Value = 1400
$ws.Range("L16").Value = 1400
$ws.Range("N16").Value = -1974

$ws.Range("H31").Value = 13694.479
$ws.Range("I31").Value = 36485.43
$ws.Range("J31").Value = 3723.4375
$ws.Range("K31").Value = 36485.43
$ws.Range("L31").Value = 3723.4375
$ws.Range("M31").Value = -36190.43
$ws.Range("N31").Value = -4313.4375

$ws.Range("H34").Value = 13694.479
$ws.Range("I34").Value = 36485.43
$ws.Range("J34").Value = 3723.4375
$ws.Range("K34").Value = 36485.43
$ws.Range("L34").Value = 3723.4375
$ws.Range("M34").Value = -36283.43
$ws.Range("N34").Value = -4127.4375

$ws.Range("H62").Value = 17080.857
$ws.Range("I62").Value = 14892.25
$ws.Range("J62").Value = 19999
$ws.Range("K62").Value = 14892.25
$ws.Range("L62").Value = 19999
$ws.Range("M62").Value = -14268.25
$ws.Range("N62").Value = -21247

$ws.Range("H65").Value = 17080.857
$ws.Range("I65").Value = 14892.25
$ws.Range("J65").Value = 19999
$ws.Range("K65").Value = 74461.25
$ws.Range("L65").Value = 99995
$ws.Range("M65").Value = -71341.25
$ws.Range("N65").Value = -106235

$ws.Range("H113").Value = 1750
$ws.Range("J113").Value = 1400
$ws.Range("L113").Value = 1400
$ws.Range("N113").Value = -5740

$ws.Range("H122").Value = 18086.4
$ws.Range("I122").Value = 1363.2727
$ws.Range("K122").Value = 4089.8181
$ws.Range("M122").Value = -1639.8181

$ws.Range("H132").Value = 1999.8462
$ws.Range("I132").Value = 1999.8462
$ws.Range("K132").Value = 5999.5386
$ws.Range("M132").Value = -3469.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6474.5654
$ws.Range("I70").Value = 6510.706
$ws.Range("J70").Value = 6372.1665
$ws.Range("K70").Value = 6510.706
$ws.Range("L70").Value = 6372.1665
$ws.Range("M70").Value = -6240.706
$ws.Range("N70").Value = -6912.1665

$ws.Range("H73").Value = 6474.5654
$ws.Range("I73").Value = 6510.706
$ws.Range("J73").Value = 6372.1665
$ws.Range("K73").Value = 6510.706
$ws.Range("L73").Value = 6372.1665
$ws.Range("M73").Value = -5574.706
$ws.Range("N73").Value = -8244.166499999999

$ws.Range("H80").Value = 7595.5835
$ws.Range("J80").Value = 8961
$ws.Range("L80").Value = 8961
$ws.Range("N80").Value = -10957

$ws.Range("H83").Value = 7595.5835
$ws.Range("J83").Value = 8961
$ws.Range("L83").Value = 44805
$ws.Range("N83").Value = -54789

$ws.Range("H97").Value = 2197.6
$ws.Range("I97").Value = 994.5
$ws.Range("K97").Value = 994.5
$ws.Range("M97").Value = -498.5

$ws.Range("H132").Value = 2530.25
$ws.Range("J132").Value = 2449.2
$ws.Range("L132").Value = 7347.599999999999
$ws.Range("N132").Value = -12407.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2500
$ws.Range("I13").Value = 2500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -2360

$ws.Range("H61").Value = 2661.6875
$ws.Range("I61").Value = 2008.091
$ws.Range("J61").Value = 4099.6
$ws.Range("K61").Value = 2008.091
$ws.Range("L61").Value = 4099.6
$ws.Range("M61").Value = -1806.091
$ws.Range("N61").Value = -4503.6

$ws.Range("H113").Value = 2661.6875
$ws.Range("I113").Value = 2008.091
$ws.Range("J113").Value = 4099.6
$ws.Range("K113").Value = 2008.091
$ws.Range("L113").Value = 4099.6
$ws.Range("M113").Value = 161.9090000000001
$ws.Range("N113").Value = -8439.6

$ws.Range("H122").Value = 3153.35
$ws.Range("I122").Value = 2813.4443
$ws.Range("J122").Value = 3431.4546
$ws.Range("K122").Value = 8440.332900000001
$ws.Range("L122").Value = 10294.3638
$ws.Range("M122").Value = -5990.332900000001
$ws.Range("N122").Value = -15194.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16212.125
$ws.Range("J45").Value = 17294.166
$ws.Range("L45").Value = 17294.166
$ws.Range("N45").Value = -18276.166

$ws.Range("H122").Value = 2893.3928
$ws.Range("I122").Value = 2972
$ws.Range("J122").Value = 2657.5715
$ws.Range("K122").Value = 8916
$ws.Range("L122").Value = 7972.7145
$ws.Range("M122").Value = -6466
$ws.Range("N122").Value = -12872.7145

$ws.Range("H132").Value = 1656.3478
$ws.Range("J132").Value = 948.5
$ws.Range("L132").Value = 2845.5
$ws.Range("N132").Value = -7905.5

$ws.Range("H136").Value = 1695.7142
$ws.Range("I136").Value = 1574.1
$ws.Range("J136").Value = 1999.75
$ws.Range("K136").Value = 4722.299999999999
$ws.Range("L136").Value = 5999.25
$ws.Range("M136").Value = -2172.299999999999
$ws.Range("N136").Value = -11099.25
